# Loan RBI, Variable Instalments
# Insert a new column before column N on the "Repayment schedule" sheet,
# select cell M13, and make "Repayment schedule" the active sheet/tab.

$wb = $excel.ActiveWorkbook

$wsRepayment = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column); formatting is
# inherited from the column to the left (M), matching native Excel behaviour.
$wsRepayment.Columns("N").Insert()
$wsRepayment.Columns("N").ColumnWidth = 10.166666666666666

# Activate the sheet (this becomes the active / selected tab of the workbook)
$wsRepayment.Activate()

# Update the selection on that sheet
$wsRepayment.Range("M13").Select()
